# JdT_nithujan.xlsx - add new journal-de-travail entries (rows 66-77)
# Row 66 is a blank separator row (same visual style as row 58/48/23/9).
# Rows 67-77 are new data rows continuing the log after row 65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- formats -------------------------------------------------------------
# Stamp the formatting for the new rows by copying it from existing rows
# that already carry the styles we need, then fill in the actual values.

# Row 66: blank separator row -> copy format from row 58 (A=12,B/C=9,D=10,E/F/G=11)
$ws.Range("A58:G58").Copy()
$ws.Range("A66:G66").PasteSpecial(-4122)

# Rows 67-77: normal data rows -> copy format from row 64 (A=7,B/C=5,D=6,E/F/G=2)
$ws.Range("A64:G64").Copy()
$ws.Range("A67:G67").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A68:G68").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A69:G69").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A70:G70").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A71:G71").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A72:G72").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A73:G73").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A74:G74").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A75:G75").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A76:G76").PasteSpecial(-4122)
$ws.Range("A64:G64").Copy()
$ws.Range("A77:G77").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Taller rows for the longer action/result text (matches the ht="30"/"90"/"45" rows)
$ws.Rows(67).RowHeight = 30
$ws.Rows(68).RowHeight = 90
$ws.Rows(72).RowHeight = 30
$ws.Rows(76).RowHeight = 45

# --- times, durations and task column (no new shared strings here) -------

$ws.Cells.Item(67, 2).Value = 0.33333333333333331
$ws.Cells.Item(67, 3).Value = 0.34375
$ws.Cells.Item(67, 4).Formula = "=C67-B67"
$ws.Cells.Item(67, 5).Value = "Implémentation"

$ws.Cells.Item(68, 2).Value = 0.34375
$ws.Cells.Item(68, 3).Value = 0.38541666666666669
$ws.Cells.Item(68, 4).Formula = "=C68-B68"
$ws.Cells.Item(68, 5).Value = "Implémentation"

$ws.Cells.Item(69, 2).Value = 0.38541666666666669
$ws.Cells.Item(69, 3).Value = 0.39930555555555558
$ws.Cells.Item(69, 4).Formula = "=C69-B69"
$ws.Cells.Item(69, 5).Value = "Implémentation"

$ws.Cells.Item(70, 2).Value = 0.40972222222222227
$ws.Cells.Item(70, 3).Value = 0.4201388888888889
$ws.Cells.Item(70, 4).Formula = "=C70-B70"
$ws.Cells.Item(70, 5).Value = "Implémentation"

$ws.Cells.Item(71, 2).Value = 0.4201388888888889
$ws.Cells.Item(71, 3).Value = 0.47916666666666669
$ws.Cells.Item(71, 4).Formula = "=C71-B71"
$ws.Cells.Item(71, 5).Value = "Implémentation"

$ws.Cells.Item(72, 2).Value = 0.47916666666666669
$ws.Cells.Item(72, 3).Value = 0.4826388888888889
$ws.Cells.Item(72, 4).Formula = "=C72-B72"
$ws.Cells.Item(72, 5).Value = "Analyse"

$ws.Cells.Item(73, 2).Value = 0.4826388888888889
$ws.Cells.Item(73, 3).Value = 0.5
$ws.Cells.Item(73, 4).Formula = "=C73-B73"
$ws.Cells.Item(73, 5).Value = "Implémentation"

$ws.Cells.Item(74, 2).Value = 0.51041666666666663
$ws.Cells.Item(74, 3).Value = 0.50347222222222221
$ws.Cells.Item(74, 4).Formula = "=C74-B74"
$ws.Cells.Item(74, 5).Value = "Documentation"

$ws.Cells.Item(75, 2).Value = 0.5625
$ws.Cells.Item(75, 3).Value = 0.62847222222222221
$ws.Cells.Item(75, 4).Formula = "=C75-B75"
$ws.Cells.Item(75, 5).Value = "Documentation"

$ws.Cells.Item(76, 2).Value = 0.63888888888888895
$ws.Cells.Item(76, 3).Value = 0.65625
$ws.Cells.Item(76, 4).Formula = "=C76-B76"
$ws.Cells.Item(76, 5).Value = "Implémentation"

$ws.Cells.Item(77, 2).Value = 0.65625
$ws.Cells.Item(77, 3).Value = 0.70486111111111116
$ws.Cells.Item(77, 4).Formula = "=C77-B77"
$ws.Cells.Item(77, 5).Value = "Documentation"

# --- action (F) / result (G) text ------------------------------------------
# Written in the same order the original author entered them (so freshly
# minted shared-string ids line up with the authoritative file): mostly
# row-by-row F-then-G, except row 67 (G before F) and G72, filled in only
# after row 73 was already written.

$ws.Cells.Item(67, 7).Value = "autre site que wikipedia utilisé: https://www.messier-objects.com/messier-catalogue/ -- ces images ont des noms qui permetttent de faciliter le renommage"
$ws.Cells.Item(67, 6).Value = "Import des images  des éléments de la bdd en local"

$ws.Cells.Item(68, 6).Value = "Création d'un script Powershell pour extraire les noms des images dans un csv"
$ws.Cells.Item(68, 7).Value = "fichiers: get-pictures-cvs.ps1, pictures.csv et images dans /assets -; src: https://docs.microsoft.com/en-us/powershell/module/microsoft.powershell.core/about/about_regular_expressions?view=powershell-7.2 -- https://stackoverflow.com/questions/54882043/powershell-variable-assignment-vs-pipeline -- https://stackoverflow.com/questions/27970441/powershell-string-does-not-contain -- "

$ws.Cells.Item(69, 6).Value = "Ajout des données du csv  dans une table temporaire"
$ws.Cells.Item(69, 7).Value = "impossible à faire erreur d'encodage "

$ws.Cells.Item(70, 6).Value = "Ajout des données du csv  dans une table temporaire"
$ws.Cells.Item(70, 7).Value = "Encodage changé mais données qui manque dans le csv --> script à changer"

$ws.Cells.Item(71, 6).Value = "Correction du script qui créé le fichier csv pour les images"
$ws.Cells.Item(71, 7).Value = "regex erroné a été changé / problème d'encodage"

$ws.Cells.Item(72, 6).Value = "Discussion avec la cdp pour questions sur mcd/mld"

$ws.Cells.Item(73, 6).Value = "Création de la bdd"
$ws.Cells.Item(73, 7).Value = "terminée + fichier backup créé"

$ws.Cells.Item(72, 7).Value = "tables temp à ne pas inclure mais mettre dans doc/réalisation -- ok de download les images depuis un autre site"

$ws.Cells.Item(74, 6).Value = "Dossier de projet"
$ws.Cells.Item(74, 7).Value = "Réalisation: import des données de la bdd --> à finir"

$ws.Cells.Item(75, 6).Value = "Dossier de projet"
$ws.Cells.Item(75, 7).Value = "Résultats des tests à rédiger"

$ws.Cells.Item(76, 6).Value = "Scripts PowerShell"
$ws.Cells.Item(76, 7).Value = "vérification et légères modif. effectuées / src: https://stackoverflow.com/questions/4724290/powershell-run-command-from-scripts-directory"

$ws.Cells.Item(77, 6).Value = "Dossier de projet"
$ws.Cells.Item(77, 7).Value = "partie implémentation à faire lire car bcp d'explications"

# --- view state --------------------------------------------------------
# Scroll the window so row 64 is at the top and select G77 (the last cell
# touched), matching the author's final cursor position.
$ws.Range("A64").Select()
$excel.ActiveWindow.ScrollRow = 64
$ws.Range("G77").Select()
